{"js": "// Update the multiplication-expression cells in the practice-sheet table.\n// Each old expression is unique in the document, so a literal (case-\n// sensitive, non-wildcard) search-and-replace per pair is safe and precise.\nconst replacements = [\n  [\"839\u00d75=\", \"876\u00d77=\"],\n  [\"407\u00d75=\", \"618\u00d77=\"],\n  [\"458\u00d77=\", \"209\u00d78=\"],\n  [\"379\u00d76=\", \"606\u00d74=\"],\n  [\"309\u00d78=\", \"106\u00d77=\"],\n  [\"291\u00d76=\", \"733\u00d78=\"],\n  [\"907\u00d73=\", \"825\u00d73=\"],\n  [\"276\u00d76=\", \"999\u00d79=\"],\n  [\"490\u00d77=\", \"628\u00d75=\"],\n  [\"366\u00d79=\", \"640\u00d75=\"],\n  [\"570\u00d76=\", \"237\u00d79=\"],\n  [\"573\u00d79=\", \"138\u00d73=\"],\n  [\"317\u00d73=\", \"501\u00d72=\"],\n  [\"951\u00d77=\", \"213\u00d79=\"],\n  [\"645\u00d73=\", \"775\u00d79=\"],\n  [\"947\u00d73=\", \"163\u00d77=\"],\n  [\"613\u00d78=\", \"819\u00d76=\"],\n  [\"847\u00d76=\", \"609\u00d77=\"],\n  [\"542\u00d75=\", \"506\u00d73=\"],\n  [\"152\u00d77=\", \"187\u00d79=\"],\n  [\"630\u00d78=\", \"864\u00d74=\"],\n  [\"201\u00d78=\", \"779\u00d73=\"],\n  [\"510\u00d74=\", \"514\u00d74=\"],\n  [\"976\u00d76=\", \"317\u00d77=\"],\n  [\"634\u00d78=\", \"335\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-expression cells in the practice-sheet table.\n# Each old expression is unique in the document, so a literal\n# (case-sensitive, non-wildcard) Find/Replace-all per pair is safe and\n# precise.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"839\u00d75=\", \"876\u00d77=\"),\n  @(\"407\u00d75=\", \"618\u00d77=\"),\n  @(\"458\u00d77=\", \"209\u00d78=\"),\n  @(\"379\u00d76=\", \"606\u00d74=\"),\n  @(\"309\u00d78=\", \"106\u00d77=\"),\n  @(\"291\u00d76=\", \"733\u00d78=\"),\n  @(\"907\u00d73=\", \"825\u00d73=\"),\n  @(\"276\u00d76=\", \"999\u00d79=\"),\n  @(\"490\u00d77=\", \"628\u00d75=\"),\n  @(\"366\u00d79=\", \"640\u00d75=\"),\n  @(\"570\u00d76=\", \"237\u00d79=\"),\n  @(\"573\u00d79=\", \"138\u00d73=\"),\n  @(\"317\u00d73=\", \"501\u00d72=\"),\n  @(\"951\u00d77=\", \"213\u00d79=\"),\n  @(\"645\u00d73=\", \"775\u00d79=\"),\n  @(\"947\u00d73=\", \"163\u00d77=\"),\n  @(\"613\u00d78=\", \"819\u00d76=\"),\n  @(\"847\u00d76=\", \"609\u00d77=\"),\n  @(\"542\u00d75=\", \"506\u00d73=\"),\n  @(\"152\u00d77=\", \"187\u00d79=\"),\n  @(\"630\u00d78=\", \"864\u00d74=\"),\n  @(\"201\u00d78=\", \"779\u00d73=\"),\n  @(\"510\u00d74=\", \"514\u00d74=\"),\n  @(\"976\u00d76=\", \"317\u00d77=\"),\n  @(\"634\u00d78=\", \"335\u00d73=\")\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $p[0]\n  $find.Replacement.Text = $p[1]\n  $find.Execute($p[0], $true, $false, $false, $false, $false, $true, 1, $false, $p[1], 2)\n}\n"}
